$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 10.2 = 41898.98 pesos`n✅ 41898.98 pesos = 10.16 = 965.65 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate cells N10/O10 and N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 98
$ws2.Range("O10").Value = 4106.1
$ws2.Range("N12").Value = 4122
$ws2.Range("O12").Value = 95
